$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 778.5143
$ws.Range("J17").Value = 778.5143
$ws.Range("L17").Value = 2335.5429
$ws.Range("N17").Value = -2671.5429

# ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 5146
$ws.Range("I69").Value = 6235
$ws.Range("J69").Value = 4420
$ws.Range("K69").Value = 18705
$ws.Range("L69").Value = 13260
$ws.Range("M69").Value = -17831
$ws.Range("N69").Value = -15008

# ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 5146
$ws.Range("I72").Value = 6235
$ws.Range("J72").Value = 4420
$ws.Range("K72").Value = 56115
$ws.Range("L72").Value = 39780
$ws.Range("M72").Value = -51747
$ws.Range("N72").Value = -48516

# ALC row 94
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 2580.4
$ws.Range("I94").Value = 1543.4286
$ws.Range("K94").Value = 1543.4286
$ws.Range("M94").Value = -1092.4286

# ALC row 96
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 2332627.8
$ws.Range("I96").Value = 3206072
$ws.Range("J96").Value = 3443
$ws.Range("K96").Value = 9618216
$ws.Range("L96").Value = 10329
$ws.Range("M96").Value = -9616843
$ws.Range("N96").Value = -13075

# ALC row 127
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H127").Value = 710.875
$ws.Range("I127").Value = 598.1429000000001
$ws.Range("J127").Value = 1500
$ws.Range("K127").Value = 1794.4287
$ws.Range("L127").Value = 4500
$ws.Range("M127").Value = 3165.5713
$ws.Range("N127").Value = -14420

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5627.431
$ws.Range("I132").Value = 4119.488
$ws.Range("J132").Value = 11810
$ws.Range("K132").Value = 12358.464
$ws.Range("L132").Value = 35430
$ws.Range("M132").Value = -9828.464
$ws.Range("N132").Value = -40490

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1746.258
$ws.Range("I137").Value = 1746.258
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5238.774
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -2688.774
$ws.Range("N137").ClearContents()

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2140.4614
$ws.Range("I138").Value = 1068.875
$ws.Range("J138").Value = 7039.143
$ws.Range("K138").Value = 3206.625
$ws.Range("L138").Value = 21117.429
$ws.Range("M138").Value = 1933.375
$ws.Range("N138").Value = -31397.429

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2457.182
$ws.Range("I2").Value = 2061.5557
$ws.Range("J2").Value = 4237.5
$ws.Range("K2").Value = 2061.5557
$ws.Range("L2").Value = 4237.5
$ws.Range("M2").Value = -1948.5557
$ws.Range("N2").Value = -4463.5

# ARM row 23
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4717.4824
$ws.Range("I32").Value = 4167.231
$ws.Range("K32").Value = 4167.231
$ws.Range("M32").Value = -3880.231

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3188.56
$ws.Range("I45").Value = 1419.2
$ws.Range("J45").Value = 4368.1333
$ws.Range("K45").Value = 1419.2
$ws.Range("L45").Value = 4368.1333
$ws.Range("M45").Value = -1042.2
$ws.Range("N45").Value = -5122.1333

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2652.2896
$ws.Range("I61").Value = 1869.4286
$ws.Range("J61").Value = 3108.9583
$ws.Range("K61").Value = 1869.4286
$ws.Range("L61").Value = 3108.9583
$ws.Range("M61").Value = -1657.4286
$ws.Range("N61").Value = -3532.9583

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2457.182
$ws.Range("I116").Value = 2061.5557
$ws.Range("J116").Value = 4237.5
$ws.Range("K116").Value = 2061.5557
$ws.Range("L116").Value = 4237.5
$ws.Range("M116").Value = 232.4443000000001
$ws.Range("N116").Value = -8825.5

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1233.7028
$ws.Range("I122").Value = 1056.7931
$ws.Range("J122").Value = 1875
$ws.Range("K122").Value = 3170.379300000001
$ws.Range("L122").Value = 5625
$ws.Range("M122").Value = -720.3793000000005
$ws.Range("N122").Value = -10525

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2652.2896
$ws.Range("I136").Value = 1869.4286
$ws.Range("J136").Value = 3108.9583
$ws.Range("K136").Value = 5608.2858
$ws.Range("L136").Value = 9326.874899999999
$ws.Range("M136").Value = -3058.2858
$ws.Range("N136").Value = -14426.8749

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2457.182
$ws.Range("I3").Value = 2061.5557
$ws.Range("J3").Value = 4237.5
$ws.Range("K3").Value = 2061.5557
$ws.Range("L3").Value = 4237.5
$ws.Range("M3").Value = -1947.5557
$ws.Range("N3").Value = -4465.5

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2408.9
$ws.Range("I105").Value = 1727.1428
$ws.Range("J105").Value = 3999.6667
$ws.Range("K105").Value = 1727.1428
$ws.Range("L105").Value = 3999.6667
$ws.Range("M105").Value = 19.85719999999992
$ws.Range("N105").Value = -7493.6667

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1161.4286
$ws.Range("I94").Value = 1200
$ws.Range("J94").Value = 1155
$ws.Range("K94").Value = 1200
$ws.Range("L94").Value = 1155
$ws.Range("M94").Value = -749
$ws.Range("N94").Value = -2057

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 378.57144
$ws.Range("I4").Value = 216.66667
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 650.00001
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -538.00001
$ws.Range("N4").Value = -1724

# CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 10206892
$ws.Range("I117").Value = 468.25
$ws.Range("J117").Value = 14289462
$ws.Range("K117").Value = 1404.75
$ws.Range("L117").Value = 42868386
$ws.Range("M117").Value = 2037.25
$ws.Range("N117").Value = -42875270

# CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 234403.4
$ws.Range("I121").Value = 296.25
$ws.Range("J121").Value = 319533.28
$ws.Range("K121").Value = 888.75
$ws.Range("L121").Value = 958599.8400000001
$ws.Range("M121").Value = 421.25
$ws.Range("N121").Value = -961219.8400000001

# GSM row 51
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 55000
$ws.Range("J51").Value = 55000
$ws.Range("L51").Value = 55000
$ws.Range("N51").Value = -56018

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2320.5757
$ws.Range("I102").Value = 1685.0869
$ws.Range("J102").Value = 3782.2
$ws.Range("K102").Value = 1685.0869
$ws.Range("L102").Value = 3782.2
$ws.Range("M102").Value = -63.08690000000001
$ws.Range("N102").Value = -7026.2

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2122.64
$ws.Range("I122").Value = 2002.75
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6008.25
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -3558.25
$ws.Range("N122").Value = -19900

# GSM row 123
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226

# LTW row 2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 4332.8335
$ws.Range("I2").Value = 3500
$ws.Range("J2").Value = 4499.4
$ws.Range("K2").Value = 3500
$ws.Range("L2").Value = 4499.4
$ws.Range("M2").Value = -3388
$ws.Range("N2").Value = -4723.4

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2751.1428
$ws.Range("I7").Value = 2071.6
$ws.Range("K7").Value = 2071.6
$ws.Range("M7").Value = -1959.6

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2751.1428
$ws.Range("I126").Value = 2071.6
$ws.Range("K126").Value = 6214.799999999999
$ws.Range("M126").Value = -3744.799999999999

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 43481656
$ws.Range("I132").Value = 50003156
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 150009468
$ws.Range("L132").Value = 14998.0005
$ws.Range("M132").Value = -150006938
$ws.Range("N132").Value = -20058.0005

# WVR row 110
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3155.111
$ws.Range("I132").Value = 3459.5
$ws.Range("J132").Value = 2976.0588
$ws.Range("K132").Value = 10378.5
$ws.Range("L132").Value = 8928.1764
$ws.Range("M132").Value = -7848.5
$ws.Range("N132").Value = -13988.1764
